$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 5).Value = 1
}

$ws.Range("E10").Select()
